$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "column1"
$ws.Range("B1").Value = "colum2"
$ws.Range("C1").Value = "column3"

# Data rows
$ws.Range("A2").Value = "v1"
$ws.Range("B3").Value = "v2"
$ws.Range("C3").Value = "v3"

# Apply thin box border around the used range A1:C3
$ws.Range("A1:C3").Borders.LineStyle = 1
$ws.Range("A1:C3").Borders.Weight = 2

# Match the saved selection position from the authored workbook
$ws.Range("N4").Select() | Out-Null
